$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextValue($cellRef, $text) {
    $target = $ws.Range($cellRef)
    $helper = $ws.Range("Z1")
    $helper.Value = "'" + $text
    $helper.Copy()
    $target.PasteSpecial(-4163)
    $helper.Clear()
}

Set-TextValue "E2" "4490"
Set-TextValue "E3" "5035"
Set-TextValue "E4" "4828"
Set-TextValue "E5" "4638"
Set-TextValue "E6" "5378"
Set-TextValue "E7" "5899"
Set-TextValue "E8" "6153"
Set-TextValue "E9" "6363"
Set-TextValue "E10" "6550"
Set-TextValue "E11" "7175"
Set-TextValue "E12" "7433"
Set-TextValue "E13" "7964"
Set-TextValue "E14" "8396"
Set-TextValue "E15" "8915"
Set-TextValue "E16" "9430"
Set-TextValue "E17" "9997"
Set-TextValue "E18" "9867"
Set-TextValue "E19" "9918"
Set-TextValue "E20" "11210"
Set-TextValue "E21" "12310"
Set-TextValue "E22" "12913"
Set-TextValue "E23" "13885"
Set-TextValue "E24" "15108"
Set-TextValue "E25" "15374"
Set-TextValue "E26" "15980"
Set-TextValue "E27" "16176"
Set-TextValue "E28" "16053"
Set-TextValue "E29" "15721"
Set-TextValue "E30" "16139"
Set-TextValue "E31" "16761"
Set-TextValue "E32" "17508"
Set-TextValue "E33" "18103"
Set-TextValue "E34" "18155"
Set-TextValue "E35" "18468"
Set-TextValue "E36" "18305"
Set-TextValue "E37" "18605"
Set-TextValue "E38" "19160"
Set-TextValue "E39" "20196"
Set-TextValue "E40" "20256"
Set-TextValue "E41" "20091"
Set-TextValue "E42" "20829"
Set-TextValue "E43" "20912.9405946244"
Set-TextValue "E44" "21531.5961416535"
Set-TextValue "E45" "21484.9793858044"
Set-TextValue "E46" "22128.8046843923"
Set-TextValue "E47" "23378.1595281805"
Set-TextValue "E48" "23979.9952749834"
Set-TextValue "E49" "24226.7311877193"
Set-TextValue "E50" "24506.1607353234"
Set-TextValue "E51" "24644.2042730325"
Set-TextValue "E52" "26058.8405870726"
Set-TextValue "E53" "25427.3193823666"
Set-TextValue "E54" "24823.552894204"
Set-TextValue "E55" "24545.2450475457"
Set-TextValue "E56" "25210.8787561371"
Set-TextValue "E57" "25686.4057017676"
Set-TextValue "E58" "26571.255443271"
Set-TextValue "E59" "27624.883102167"
Set-TextValue "E60" "27857.6186676969"
Set-TextValue "E61" "27644.6974945809"
Set-TextValue "E62" "28575.531463221"

# Append new rows 63-68 for years 2011-2016
$ws.Range("A63").Value = 376
$ws.Range("B63").Value = "Israel"
$ws.Range("C63").Value = "GDP per Capita"
$ws.Range("D63").Value = 2011
Set-TextValue "E63" "29370"

$ws.Range("A64").Value = 376
$ws.Range("B64").Value = "Israel"
$ws.Range("C64").Value = "GDP per Capita"
$ws.Range("D64").Value = 2012
Set-TextValue "E64" "29517"

$ws.Range("A65").Value = 376
$ws.Range("B65").Value = "Israel"
$ws.Range("C65").Value = "GDP per Capita"
$ws.Range("D65").Value = 2013
Set-TextValue "E65" "30257"

$ws.Range("A66").Value = 376
$ws.Range("B66").Value = "Israel"
$ws.Range("C66").Value = "GDP per Capita"
$ws.Range("D66").Value = 2014
Set-TextValue "E66" "30686"

$ws.Range("A67").Value = 376
$ws.Range("B67").Value = "Israel"
$ws.Range("C67").Value = "GDP per Capita"
$ws.Range("D67").Value = 2015
Set-TextValue "E67" "30954"

$ws.Range("A68").Value = 376
$ws.Range("B68").Value = "Israel"
$ws.Range("C68").Value = "GDP per Capita"
$ws.Range("D68").Value = 2016
Set-TextValue "E68" "31701"

Write-Host "Update complete"
